# Generate Report for Handoff
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#   for the entries that just got a new handoff xliff generated (06c2d3c9-... file).
# - Update its Priority from "low" to "ht" on both locale sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for all data rows
$wsOverview.Range("G4").Value = "2016-09-06 10:40:29"
$wsOverview.Range("G5").Value = "2016-09-06 10:40:29"
$wsOverview.Range("G6").Value = "2016-09-06 10:40:29"
$wsOverview.Range("G7").Value = "2016-09-06 10:40:29"

# zh-cn sheet: Priority (E) low -> ht, Latest Handoff Datetime (H) refreshed
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("E7").Value = "ht"

$wsZhCn.Range("H4").Value = "2016-09-06 10:40:22"
$wsZhCn.Range("H5").Value = "2016-09-06 10:40:22"
$wsZhCn.Range("H6").Value = "2016-09-06 10:40:22"
$wsZhCn.Range("H7").Value = "2016-09-06 10:40:22"

# de-de sheet: Priority (E) low -> ht, Latest Handoff Datetime (H) refreshed
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("E7").Value = "ht"

$wsDeDe.Range("H4").Value = "2016-09-06 10:40:29"
$wsDeDe.Range("H5").Value = "2016-09-06 10:40:29"
$wsDeDe.Range("H6").Value = "2016-09-06 10:40:29"
$wsDeDe.Range("H7").Value = "2016-09-06 10:40:29"
